# Additional companies sent for questionaire
# Remove the "Parent company" and "Location County/City" columns from the
# locomotive list (these fields are no longer part of the collected data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Parent company"). Remaining columns shift left by one,
# so the former column E ("Location County/City") becomes column D.
$ws.Range("B1").EntireColumn.Delete()

# Delete the (now) column D ("Location County/City").
$ws.Range("D1").EntireColumn.Delete()

# Reflect the final cell selection left behind after the edit.
$ws.Range("L8").Select()
